$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# --- Row 15: previously-blank cells now carry the literal text "nan" ---
# (L15, N15, O15 already hold real values and are left untouched)
$ws.Range("B15:K15").Value = "nan"
$ws.Range("M15").Value = "nan"
$ws.Range("P15").Value = "nan"

# --- Row 16: brand-new maintenance event row ---
# A16 holds the text "24" (the card id, same text used throughout column A),
# so force Text format before assigning it - otherwise Excel would read
# "24" as a number instead of a string.
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "24"
$ws.Range("A16").Style = "Normal"

# The numeric/boolean columns (B:K) and the last "Serviced by" column (P)
# stay blank for this event, but the cells themselves still get materialised
# (matching the rest of the sheet, where every row has a full A:P cell set).
$ws.Range("B16:K16").NumberFormat = "@"
$ws.Range("B16:K16").Style = "Normal"
$ws.Range("P16").NumberFormat = "@"
$ws.Range("P16").Style = "Normal"

$ws.Range("L16").Value = "20\5\2025"
$ws.Range("M16").Value = "632.3 t"
$ws.Range("N16").Value = "تم عمل صيانه وسن السلندر وتغير الجرائد الاماميه(1_2_4_5_7_8) وتغير الجرائد الخلفيه(1_5_8)"
$ws.Range("O16").Value = "الخبير"
